# Auto-generated edit script: applies updated market-data snapshot values
# to the Leve profit columns (H-N) across all 8 job sheets, per the scheduled
# runner's commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 31
$ws.Range("H31").Value = 15071.625
$ws.Range("I31").Value = 15071.625
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 45214.875
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -44984.875
$ws.Range("N31").ClearContents()
# Row 62
$ws.Range("H62").Value = 1459.8572
$ws.Range("I62").Value = 1572.4667
$ws.Range("J62").Value = 1178.3334
$ws.Range("K62").Value = 1572.4667
$ws.Range("L62").Value = 1178.3334
$ws.Range("M62").Value = -948.4666999999999
$ws.Range("N62").Value = -2426.3334
# Row 65
$ws.Range("H65").Value = 1459.8572
$ws.Range("I65").Value = 1572.4667
$ws.Range("J65").Value = 1178.3334
$ws.Range("K65").Value = 7862.3335
$ws.Range("L65").Value = 5891.666999999999
$ws.Range("M65").Value = -4742.3335
$ws.Range("N65").Value = -12131.667
# Row 106
$ws.Range("H106").Value = 2552.5
$ws.Range("I106").Value = 2815.625
$ws.Range("J106").Value = 1500
$ws.Range("K106").Value = 2815.625
$ws.Range("L106").Value = 1500
$ws.Range("M106").Value = -2184.625
$ws.Range("N106").Value = -2762
# Row 137
$ws.Range("H137").Value = 1050.8695
$ws.Range("I137").Value = 1022.3415
$ws.Range("J137").Value = 1284.8
$ws.Range("K137").Value = 3067.0245
$ws.Range("L137").Value = 3854.4
$ws.Range("M137").Value = -517.0245
$ws.Range("N137").Value = -8954.4

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 12506.623
$ws.Range("I32").Value = 14622.816
$ws.Range("J32").Value = 3865.5
$ws.Range("K32").Value = 14622.816
$ws.Range("L32").Value = 3865.5
$ws.Range("M32").Value = -14335.816
$ws.Range("N32").Value = -4439.5
# Row 61
$ws.Range("H61").Value = 1676.5581
$ws.Range("I61").Value = 1520.8422
$ws.Range("J61").Value = 2860
$ws.Range("K61").Value = 1520.8422
$ws.Range("L61").Value = 2860
$ws.Range("M61").Value = -1308.8422
$ws.Range("N61").Value = -3284
# Row 74
$ws.Range("H74").Value = 871.7143
$ws.Range("I74").Value = 684.5
$ws.Range("J74").Value = 1995
$ws.Range("K74").Value = 684.5
$ws.Range("L74").Value = 1995
$ws.Range("M74").Value = 189.5
$ws.Range("N74").Value = -3743
# Row 77
$ws.Range("H77").Value = 871.7143
$ws.Range("I77").Value = 684.5
$ws.Range("J77").Value = 1995
$ws.Range("K77").Value = 3422.5
$ws.Range("L77").Value = 9975
$ws.Range("M77").Value = 945.5
$ws.Range("N77").Value = -18711
# Row 97
$ws.Range("H97").Value = 967.5
$ws.Range("I97").Value = 940
$ws.Range("J97").Value = 1013.3333
$ws.Range("K97").Value = 940
$ws.Range("L97").Value = 1013.3333
$ws.Range("M97").Value = -444
$ws.Range("N97").Value = -2005.3333
# Row 122
$ws.Range("H122").Value = 3509.9285
$ws.Range("I122").Value = 3016.5
$ws.Range("J122").Value = 4167.8335
$ws.Range("K122").Value = 9049.5
$ws.Range("L122").Value = 12503.5005
$ws.Range("M122").Value = -6599.5
$ws.Range("N122").Value = -17403.5005
# Row 123
$ws.Range("H123").Value = 24225
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 24225
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 24225
$ws.Range("N123").Value = -34025
# Row 131
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
# Row 132
$ws.Range("H132").Value = 1596.0615
$ws.Range("I132").Value = 1038.238
$ws.Range("J132").Value = 2614.6956
$ws.Range("K132").Value = 3114.714
$ws.Range("L132").Value = 7844.0868
$ws.Range("M132").Value = -584.7139999999999
$ws.Range("N132").Value = -12904.0868
# Row 136
$ws.Range("H136").Value = 1676.5581
$ws.Range("I136").Value = 1520.8422
$ws.Range("J136").Value = 2860
$ws.Range("K136").Value = 4562.5266
$ws.Range("L136").Value = 8580
$ws.Range("M136").Value = -2012.5266
$ws.Range("N136").Value = -13680

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 1915.2927
$ws.Range("I134").Value = 1386.85
$ws.Range("J134").Value = 2418.5715
$ws.Range("K134").Value = 4160.549999999999
$ws.Range("L134").Value = 7255.7145
$ws.Range("M134").Value = -1625.549999999999
$ws.Range("N134").Value = -12325.7145

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1762.8334
$ws.Range("I31").Value = 1475
$ws.Range("J31").Value = 2683.9
$ws.Range("K31").Value = 1475
$ws.Range("L31").Value = 2683.9
$ws.Range("M31").Value = -1180
$ws.Range("N31").Value = -3273.9
# Row 34
$ws.Range("H34").Value = 1762.8334
$ws.Range("I34").Value = 1475
$ws.Range("J34").Value = 2683.9
$ws.Range("K34").Value = 1475
$ws.Range("L34").Value = 2683.9
$ws.Range("M34").Value = -1273
$ws.Range("N34").Value = -3087.9
# Row 58
$ws.Range("H58").Value = 951318.4399999999
$ws.Range("I58").Value = 1684877.9
$ws.Range("J58").Value = 2006.2354
$ws.Range("K58").Value = 1684877.9
$ws.Range("L58").Value = 2006.2354
$ws.Range("M58").Value = -1684674.9
$ws.Range("N58").Value = -2412.2354
# Row 107
$ws.Range("H107").Value = 330.5
$ws.Range("I107").Value = 317.625
$ws.Range("J107").Value = 356.25
$ws.Range("K107").Value = 317.625
$ws.Range("L107").Value = 356.25
$ws.Range("M107").Value = 1602.375
$ws.Range("N107").Value = -4196.25
# Row 132
$ws.Range("H132").Value = 411292.06
$ws.Range("I132").Value = 484005
$ws.Range("J132").Value = 4099.6
$ws.Range("K132").Value = 1452015
$ws.Range("L132").Value = 12298.8
$ws.Range("M132").Value = -1449485
$ws.Range("N132").Value = -17358.8
# Row 134
$ws.Range("H134").Value = 1671.9803
$ws.Range("I134").Value = 1269.1892
$ws.Range("J134").Value = 2736.5
$ws.Range("K134").Value = 3807.5676
$ws.Range("L134").Value = 8209.5
$ws.Range("M134").Value = -1272.5676
$ws.Range("N134").Value = -13279.5
# Row 136
$ws.Range("H136").Value = 951318.4399999999
$ws.Range("I136").Value = 1684877.9
$ws.Range("J136").Value = 2006.2354
$ws.Range("K136").Value = 5054633.699999999
$ws.Range("L136").Value = 6018.706200000001
$ws.Range("M136").Value = -5052083.699999999
$ws.Range("N136").Value = -11118.7062

$ws = $wb.Worksheets.Item("CUL")
# Row 132
$ws.Range("H132").Value = 1676.762
$ws.Range("I132").Value = 1167.3334
$ws.Range("J132").Value = 1880.5333
$ws.Range("K132").Value = 10506.0006
$ws.Range("L132").Value = 16924.7997
$ws.Range("M132").Value = -7976.000599999999
$ws.Range("N132").Value = -21984.7997

$ws = $wb.Worksheets.Item("GSM")
# Row 82
$ws.Range("H82").Value = 31200
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 31200
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 31200
$ws.Range("N82").Value = -31966
# Row 85
$ws.Range("H85").Value = 31200
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 31200
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 31200
$ws.Range("N85").Value = -33852
# Row 93
$ws.Range("H93").Value = 40000
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 40000
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 40000
$ws.Range("N93").Value = -43744
# Row 97
$ws.Range("H97").Value = 53426.6
$ws.Range("I97").Value = 66033.25
$ws.Range("J97").Value = 3000
$ws.Range("K97").Value = 66033.25
$ws.Range("L97").Value = 3000
$ws.Range("M97").Value = -65537.25
$ws.Range("N97").Value = -3992
# Row 113
$ws.Range("H113").Value = 1497.3334
$ws.Range("I113").Value = 764.2
$ws.Range("J113").Value = 1863.9
$ws.Range("K113").Value = 764.2
$ws.Range("L113").Value = 1863.9
$ws.Range("M113").Value = 1405.8
$ws.Range("N113").Value = -6203.9
# Row 122
$ws.Range("H122").Value = 2904.4814
$ws.Range("I122").Value = 2031.3077
$ws.Range("J122").Value = 3715.2856
$ws.Range("K122").Value = 6093.9231
$ws.Range("L122").Value = 11145.8568
$ws.Range("M122").Value = -3643.9231
$ws.Range("N122").Value = -16045.8568
# Row 123
$ws.Range("H123").Value = 17347.334
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 17347.334
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 17347.334
$ws.Range("N123").Value = -22247.334
# Row 132
$ws.Range("H132").Value = 2733
$ws.Range("I132").Value = 2084.8667
$ws.Range("J132").Value = 3948.25
$ws.Range("K132").Value = 6254.6001
$ws.Range("L132").Value = 11844.75
$ws.Range("M132").Value = -3724.6001
$ws.Range("N132").Value = -16904.75

$ws = $wb.Worksheets.Item("LTW")
# Row 38
$ws.Range("H38").Value = 35000
$ws.Range("I38").Value = 50000
$ws.Range("J38").Value = 20000
$ws.Range("K38").Value = 50000
$ws.Range("L38").Value = 20000
$ws.Range("M38").Value = -49590
$ws.Range("N38").Value = -20820
# Row 93
$ws.Range("H93").Value = 1341.4667
$ws.Range("I93").Value = 725
$ws.Range("J93").Value = 1565.6364
$ws.Range("K93").Value = 725
$ws.Range("L93").Value = 1565.6364
$ws.Range("M93").Value = 523
$ws.Range("N93").Value = -4061.6364
# Row 122
$ws.Range("H122").Value = 18187620
$ws.Range("I122").Value = 6258.2856
$ws.Range("J122").Value = 50005000
$ws.Range("K122").Value = 18774.8568
$ws.Range("L122").Value = 150015000
$ws.Range("M122").Value = -16324.8568
$ws.Range("N122").Value = -150019900
# Row 132
$ws.Range("H132").Value = 4557
$ws.Range("I132").Value = 4617.4
$ws.Range("J132").Value = 4389.222
$ws.Range("K132").Value = 13852.2
$ws.Range("L132").Value = 13167.666
$ws.Range("M132").Value = -11322.2
$ws.Range("N132").Value = -18227.666
# Row 136
$ws.Range("H136").Value = 27026172
$ws.Range("I136").Value = 41668200
$ws.Range("J136").Value = 670519.3
$ws.Range("K136").Value = 125004600
$ws.Range("L136").Value = 2011557.9
$ws.Range("M136").Value = -125002050
$ws.Range("N136").Value = -2016657.9

$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 1643.2858
$ws.Range("I96").Value = 1301.5
$ws.Range("J96").Value = 1780
$ws.Range("K96").Value = 1301.5
$ws.Range("L96").Value = 1780
$ws.Range("M96").Value = 71.5
$ws.Range("N96").Value = -4526
# Row 122
$ws.Range("H122").Value = 41670228
$ws.Range("I122").Value = 66669670
$ws.Range("J122").Value = 4486.778
$ws.Range("K122").Value = 200009010
$ws.Range("L122").Value = 13460.334
$ws.Range("M122").Value = -200006560
$ws.Range("N122").Value = -18360.334
# Row 123
$ws.Range("H123").Value = 37875.375
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 37875.375
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 37875.375
$ws.Range("N123").Value = -47675.375
# Row 132
$ws.Range("H132").Value = 2640.0527
$ws.Range("I132").Value = 1996.1428
$ws.Range("J132").Value = 3015.6667
$ws.Range("K132").Value = 5988.428400000001
$ws.Range("L132").Value = 9047.000100000001
$ws.Range("M132").Value = -3458.428400000001
$ws.Range("N132").Value = -14107.0001
# Row 136
$ws.Range("H136").Value = 1228.6666
$ws.Range("I136").Value = 989.5789
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 2968.7367
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -418.7366999999999
$ws.Range("N136").Value = -15600
